# Generate Report for Handback
# Refresh the handoff/handback timestamps for the second tracked file
# (d5289b5b-7386-4bf1-ae40-3302f0f8189c.md) now that its localization
# round-trip has completed.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: bump the "Latest HO Xliff Generate Date" for row 3
# (d5289b5b-...) to reflect the newly generated de-de handback xliff.
$overview.Range("G3").Value = "2016-09-06 05:01:51"

# zh-cn sheet: row 3 picked up a new handoff/handback cycle.
$zhcn.Range("H3").Value = "2016-09-06 05:01:46"
$zhcn.Range("K3").Value = "2016-09-06 05:02:10"

# de-de sheet: row 3 picked up a new handoff/handback cycle.
$dede.Range("H3").Value = "2016-09-06 05:01:51"
$dede.Range("K3").Value = "2016-09-06 05:02:19"
